$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 237 <- source row 238
$ws.Cells.Item(237, 2).Value = 3899309
$ws.Cells.Item(237, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(237, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(237, 5).Value = 44535.60416666666
$ws.Cells.Item(237, 6).Value = 'AEK Athens'
$ws.Cells.Item(237, 7).Value = 'Panathinaikos'
$ws.Cells.Item(237, 8).Value = 1
$ws.Cells.Item(237, 9).Value = 0
$ws.Cells.Item(237, 10).Value = 'H'
$ws.Cells.Item(237, 11).Value = 1.833
$ws.Cells.Item(237, 12).Value = 3.2
$ws.Cells.Item(237, 13).Value = 4.75
$ws.Cells.Item(237, 14).Value = 1.85
$ws.Cells.Item(237, 15).Value = 3.2
$ws.Cells.Item(237, 16).Value = 5
$ws.Cells.Item(237, 17).Value = -0.5
$ws.Cells.Item(237, 18).Value = 1.8
$ws.Cells.Item(237, 19).Value = 2.05
$ws.Cells.Item(237, 20).Value = 2.25
$ws.Cells.Item(237, 21).Value = 2
$ws.Cells.Item(237, 22).Value = 1.85
$ws.Cells.Item(237, 23).Value = 0.8500000000000001
$ws.Cells.Item(237, 24).Value = -1
$ws.Cells.Item(237, 25).Value = -1
$ws.Cells.Item(237, 26).Value = 0.8
$ws.Cells.Item(237, 27).Value = -1
$ws.Cells.Item(237, 28).Value = -1
$ws.Cells.Item(237, 29).Value = 0.8500000000000001

# Row 238 <- source row 237
$ws.Cells.Item(238, 2).Value = 3898702
$ws.Cells.Item(238, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(238, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(238, 5).Value = 44535.60416666666
$ws.Cells.Item(238, 6).Value = 'Aris Salonika'
$ws.Cells.Item(238, 7).Value = 'Asteras Tripolis'
$ws.Cells.Item(238, 8).Value = 1
$ws.Cells.Item(238, 9).Value = 0
$ws.Cells.Item(238, 10).Value = 'H'
$ws.Cells.Item(238, 11).Value = 1.65
$ws.Cells.Item(238, 12).Value = 3.4
$ws.Cells.Item(238, 13).Value = 6
$ws.Cells.Item(238, 14).Value = 1.727
$ws.Cells.Item(238, 15).Value = 3.4
$ws.Cells.Item(238, 16).Value = 5.75
$ws.Cells.Item(238, 17).Value = -0.75
$ws.Cells.Item(238, 18).Value = 1.975
$ws.Cells.Item(238, 19).Value = 1.875
$ws.Cells.Item(238, 20).Value = 2
$ws.Cells.Item(238, 21).Value = 1.875
$ws.Cells.Item(238, 22).Value = 1.975
$ws.Cells.Item(238, 23).Value = 0.7270000000000001
$ws.Cells.Item(238, 24).Value = -1
$ws.Cells.Item(238, 25).Value = -1
$ws.Cells.Item(238, 26).Value = 0.4875
$ws.Cells.Item(238, 27).Value = -0.5
$ws.Cells.Item(238, 28).Value = -1
$ws.Cells.Item(238, 29).Value = 0.9750000000000001

# Row 282 <- source row 283
$ws.Cells.Item(282, 2).Value = 3898756
$ws.Cells.Item(282, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(282, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(282, 5).Value = 44590.41666666666
$ws.Cells.Item(282, 6).Value = 'Lamia'
$ws.Cells.Item(282, 7).Value = 'OFI Crete'
$ws.Cells.Item(282, 8).Value = 2
$ws.Cells.Item(282, 9).Value = 1
$ws.Cells.Item(282, 10).Value = 'H'
$ws.Cells.Item(282, 11).Value = 2.625
$ws.Cells.Item(282, 12).Value = 3.1
$ws.Cells.Item(282, 13).Value = 2.7
$ws.Cells.Item(282, 14).Value = 2.6
$ws.Cells.Item(282, 15).Value = 2.875
$ws.Cells.Item(282, 16).Value = 3.1
$ws.Cells.Item(282, 17).Value = -0.25
$ws.Cells.Item(282, 18).Value = 2.1
$ws.Cells.Item(282, 19).Value = 1.775
$ws.Cells.Item(282, 20).Value = 1.75
$ws.Cells.Item(282, 21).Value = 1.85
$ws.Cells.Item(282, 22).Value = 2
$ws.Cells.Item(282, 23).Value = 1.6
$ws.Cells.Item(282, 24).Value = -1
$ws.Cells.Item(282, 25).Value = -1
$ws.Cells.Item(282, 26).Value = 1.1
$ws.Cells.Item(282, 27).Value = -1
$ws.Cells.Item(282, 28).Value = 0.8500000000000001
$ws.Cells.Item(282, 29).Value = -1

# Row 283 <- source row 282
$ws.Cells.Item(283, 2).Value = 3898757
$ws.Cells.Item(283, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(283, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(283, 5).Value = 44590.41666666666
$ws.Cells.Item(283, 6).Value = 'Panathinaikos'
$ws.Cells.Item(283, 7).Value = 'Asteras Tripolis'
$ws.Cells.Item(283, 8).Value = 0
$ws.Cells.Item(283, 9).Value = 1
$ws.Cells.Item(283, 10).Value = 'A'
$ws.Cells.Item(283, 11).Value = 1.833
$ws.Cells.Item(283, 12).Value = 3.2
$ws.Cells.Item(283, 13).Value = 4.75
$ws.Cells.Item(283, 14).Value = 1.666
$ws.Cells.Item(283, 15).Value = 3.4
$ws.Cells.Item(283, 16).Value = 6
$ws.Cells.Item(283, 17).Value = -0.75
$ws.Cells.Item(283, 18).Value = 1.9
$ws.Cells.Item(283, 19).Value = 1.95
$ws.Cells.Item(283, 20).Value = 2
$ws.Cells.Item(283, 21).Value = 1.975
$ws.Cells.Item(283, 22).Value = 1.875
$ws.Cells.Item(283, 23).Value = -1
$ws.Cells.Item(283, 24).Value = -1
$ws.Cells.Item(283, 25).Value = 5
$ws.Cells.Item(283, 26).Value = -1
$ws.Cells.Item(283, 27).Value = 0.95
$ws.Cells.Item(283, 28).Value = -1
$ws.Cells.Item(283, 29).Value = 0.875

# Row 290 <- source row 291
$ws.Cells.Item(290, 2).Value = 3898763
$ws.Cells.Item(290, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(290, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(290, 5).Value = 44594.51041666666
$ws.Cells.Item(290, 6).Value = 'Olympiakos'
$ws.Cells.Item(290, 7).Value = 'Panetolikos'
$ws.Cells.Item(290, 8).Value = 3
$ws.Cells.Item(290, 9).Value = 1
$ws.Cells.Item(290, 10).Value = 'H'
$ws.Cells.Item(290, 11).Value = 1.25
$ws.Cells.Item(290, 12).Value = 5.75
$ws.Cells.Item(290, 13).Value = 9
$ws.Cells.Item(290, 14).Value = 1.363
$ws.Cells.Item(290, 15).Value = 5.25
$ws.Cells.Item(290, 16).Value = 7.5
$ws.Cells.Item(290, 17).Value = -1.5
$ws.Cells.Item(290, 18).Value = 2.05
$ws.Cells.Item(290, 19).Value = 1.8
$ws.Cells.Item(290, 20).Value = 2.75
$ws.Cells.Item(290, 21).Value = 1.975
$ws.Cells.Item(290, 22).Value = 1.875
$ws.Cells.Item(290, 23).Value = 0.363
$ws.Cells.Item(290, 24).Value = -1
$ws.Cells.Item(290, 25).Value = -1
$ws.Cells.Item(290, 26).Value = 1.05
$ws.Cells.Item(290, 27).Value = -1
$ws.Cells.Item(290, 28).Value = 0.9750000000000001
$ws.Cells.Item(290, 29).Value = -1

# Row 291 <- source row 290
$ws.Cells.Item(291, 2).Value = 3898762
$ws.Cells.Item(291, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(291, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(291, 5).Value = 44594.51041666666
$ws.Cells.Item(291, 6).Value = 'OFI Crete'
$ws.Cells.Item(291, 7).Value = 'Panathinaikos'
$ws.Cells.Item(291, 8).Value = 3
$ws.Cells.Item(291, 9).Value = 2
$ws.Cells.Item(291, 10).Value = 'H'
$ws.Cells.Item(291, 11).Value = 3.1
$ws.Cells.Item(291, 12).Value = 3
$ws.Cells.Item(291, 13).Value = 2.3
$ws.Cells.Item(291, 14).Value = 3.8
$ws.Cells.Item(291, 15).Value = 3.2
$ws.Cells.Item(291, 16).Value = 2.1
$ws.Cells.Item(291, 17).Value = 0.25
$ws.Cells.Item(291, 18).Value = 2.1
$ws.Cells.Item(291, 19).Value = 1.775
$ws.Cells.Item(291, 20).Value = 2
$ws.Cells.Item(291, 21).Value = 1.875
$ws.Cells.Item(291, 22).Value = 1.975
$ws.Cells.Item(291, 23).Value = 2.8
$ws.Cells.Item(291, 24).Value = -1
$ws.Cells.Item(291, 25).Value = -1
$ws.Cells.Item(291, 26).Value = 1.1
$ws.Cells.Item(291, 27).Value = -1
$ws.Cells.Item(291, 28).Value = 0.875
$ws.Cells.Item(291, 29).Value = -1

# Row 292 <- source row 293
$ws.Cells.Item(292, 2).Value = 3899314
$ws.Cells.Item(292, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(292, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(292, 5).Value = 44594.60416666666
$ws.Cells.Item(292, 6).Value = 'Aris Salonika'
$ws.Cells.Item(292, 7).Value = 'AEK Athens'
$ws.Cells.Item(292, 8).Value = 2
$ws.Cells.Item(292, 9).Value = 1
$ws.Cells.Item(292, 10).Value = 'H'
$ws.Cells.Item(292, 11).Value = 3.3
$ws.Cells.Item(292, 12).Value = 3
$ws.Cells.Item(292, 13).Value = 2.2
$ws.Cells.Item(292, 14).Value = 2.4
$ws.Cells.Item(292, 15).Value = 2.9
$ws.Cells.Item(292, 16).Value = 3.3
$ws.Cells.Item(292, 17).Value = -0.25
$ws.Cells.Item(292, 18).Value = 2.025
$ws.Cells.Item(292, 19).Value = 1.825
$ws.Cells.Item(292, 20).Value = 1.75
$ws.Cells.Item(292, 21).Value = 1.8
$ws.Cells.Item(292, 22).Value = 2.05
$ws.Cells.Item(292, 23).Value = 1.4
$ws.Cells.Item(292, 24).Value = -1
$ws.Cells.Item(292, 25).Value = -1
$ws.Cells.Item(292, 26).Value = 1.025
$ws.Cells.Item(292, 27).Value = -1
$ws.Cells.Item(292, 28).Value = 0.8
$ws.Cells.Item(292, 29).Value = -1

# Row 293 <- source row 292
$ws.Cells.Item(293, 2).Value = 3898765
$ws.Cells.Item(293, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(293, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(293, 5).Value = 44594.60416666666
$ws.Cells.Item(293, 6).Value = 'Volos NFC'
$ws.Cells.Item(293, 7).Value = 'Ionikos Nikea'
$ws.Cells.Item(293, 8).Value = 1
$ws.Cells.Item(293, 9).Value = 1
$ws.Cells.Item(293, 10).Value = 'D'
$ws.Cells.Item(293, 11).Value = 1.8
$ws.Cells.Item(293, 12).Value = 3.3
$ws.Cells.Item(293, 13).Value = 4.333
$ws.Cells.Item(293, 14).Value = 2.15
$ws.Cells.Item(293, 15).Value = 3.4
$ws.Cells.Item(293, 16).Value = 3.4
$ws.Cells.Item(293, 17).Value = -0.25
$ws.Cells.Item(293, 18).Value = 1.85
$ws.Cells.Item(293, 19).Value = 2
$ws.Cells.Item(293, 20).Value = 2.25
$ws.Cells.Item(293, 21).Value = 1.925
$ws.Cells.Item(293, 22).Value = 1.925
$ws.Cells.Item(293, 23).Value = -1
$ws.Cells.Item(293, 24).Value = 2.4
$ws.Cells.Item(293, 25).Value = -1
$ws.Cells.Item(293, 26).Value = -0.5
$ws.Cells.Item(293, 27).Value = 0.5
$ws.Cells.Item(293, 28).Value = -0.5
$ws.Cells.Item(293, 29).Value = 0.4625

# Row 316 <- source row 317
$ws.Cells.Item(316, 2).Value = 3899318
$ws.Cells.Item(316, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(316, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(316, 5).Value = 44612.60416666666
$ws.Cells.Item(316, 6).Value = 'Aris Salonika'
$ws.Cells.Item(316, 7).Value = 'PAOK Salonika'
$ws.Cells.Item(316, 8).Value = 0
$ws.Cells.Item(316, 9).Value = 0
$ws.Cells.Item(316, 10).Value = 'D'
$ws.Cells.Item(316, 11).Value = 2.875
$ws.Cells.Item(316, 12).Value = 3.2
$ws.Cells.Item(316, 13).Value = 2.5
$ws.Cells.Item(316, 14).Value = 2.7
$ws.Cells.Item(316, 15).Value = 2.875
$ws.Cells.Item(316, 16).Value = 3
$ws.Cells.Item(316, 17).Value = 0
$ws.Cells.Item(316, 18).Value = 1.8
$ws.Cells.Item(316, 19).Value = 2.05
$ws.Cells.Item(316, 20).Value = 1.75
$ws.Cells.Item(316, 21).Value = 1.85
$ws.Cells.Item(316, 22).Value = 2
$ws.Cells.Item(316, 23).Value = -1
$ws.Cells.Item(316, 24).Value = 1.875
$ws.Cells.Item(316, 25).Value = -1
$ws.Cells.Item(316, 26).Value = 0
$ws.Cells.Item(316, 27).Value = -0
$ws.Cells.Item(316, 28).Value = -1
$ws.Cells.Item(316, 29).Value = 1

# Row 317 <- source row 316
$ws.Cells.Item(317, 2).Value = 3899317
$ws.Cells.Item(317, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(317, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(317, 5).Value = 44612.60416666666
$ws.Cells.Item(317, 6).Value = 'AEK Athens'
$ws.Cells.Item(317, 7).Value = 'Giannina'
$ws.Cells.Item(317, 8).Value = 2
$ws.Cells.Item(317, 9).Value = 0
$ws.Cells.Item(317, 10).Value = 'H'
$ws.Cells.Item(317, 11).Value = 1.666
$ws.Cells.Item(317, 12).Value = 3.75
$ws.Cells.Item(317, 13).Value = 5
$ws.Cells.Item(317, 14).Value = 1.333
$ws.Cells.Item(317, 15).Value = 5
$ws.Cells.Item(317, 16).Value = 9.5
$ws.Cells.Item(317, 17).Value = -1.25
$ws.Cells.Item(317, 18).Value = 2.025
$ws.Cells.Item(317, 19).Value = 1.825
$ws.Cells.Item(317, 20).Value = 2.25
$ws.Cells.Item(317, 21).Value = 1.825
$ws.Cells.Item(317, 22).Value = 2.025
$ws.Cells.Item(317, 23).Value = 0.333
$ws.Cells.Item(317, 24).Value = -1
$ws.Cells.Item(317, 25).Value = -1
$ws.Cells.Item(317, 26).Value = 1.025
$ws.Cells.Item(317, 27).Value = -1
$ws.Cells.Item(317, 28).Value = -0.5
$ws.Cells.Item(317, 29).Value = 0.5125

# Row 421 <- source row 422
$ws.Cells.Item(421, 2).Value = 5374129
$ws.Cells.Item(421, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(421, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(421, 5).Value = 44815.64583333334
$ws.Cells.Item(421, 6).Value = 'Panathinaikos'
$ws.Cells.Item(421, 7).Value = 'AEK Athens'
$ws.Cells.Item(421, 8).Value = 2
$ws.Cells.Item(421, 9).Value = 1
$ws.Cells.Item(421, 10).Value = 'H'
$ws.Cells.Item(421, 11).Value = 2.25
$ws.Cells.Item(421, 12).Value = 3.25
$ws.Cells.Item(421, 13).Value = 3.25
$ws.Cells.Item(421, 14).Value = 2.25
$ws.Cells.Item(421, 15).Value = 3.2
$ws.Cells.Item(421, 16).Value = 3.4
$ws.Cells.Item(421, 17).Value = -0.25
$ws.Cells.Item(421, 18).Value = 1.975
$ws.Cells.Item(421, 19).Value = 1.875
$ws.Cells.Item(421, 20).Value = 2
$ws.Cells.Item(421, 21).Value = 1.825
$ws.Cells.Item(421, 22).Value = 2.025
$ws.Cells.Item(421, 23).Value = 1.25
$ws.Cells.Item(421, 24).Value = -1
$ws.Cells.Item(421, 25).Value = -1
$ws.Cells.Item(421, 26).Value = 0.9750000000000001
$ws.Cells.Item(421, 27).Value = -1
$ws.Cells.Item(421, 28).Value = 0.825
$ws.Cells.Item(421, 29).Value = -1

# Row 422 <- source row 421
$ws.Cells.Item(422, 2).Value = 5374127
$ws.Cells.Item(422, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(422, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(422, 5).Value = 44815.64583333334
$ws.Cells.Item(422, 6).Value = 'Asteras Tripolis'
$ws.Cells.Item(422, 7).Value = 'Aris Salonika'
$ws.Cells.Item(422, 8).Value = 0
$ws.Cells.Item(422, 9).Value = 2
$ws.Cells.Item(422, 10).Value = 'A'
$ws.Cells.Item(422, 11).Value = 2.6
$ws.Cells.Item(422, 12).Value = 3.1
$ws.Cells.Item(422, 13).Value = 2.7
$ws.Cells.Item(422, 14).Value = 3.5
$ws.Cells.Item(422, 15).Value = 3.2
$ws.Cells.Item(422, 16).Value = 2.2
$ws.Cells.Item(422, 17).Value = 0.25
$ws.Cells.Item(422, 18).Value = 2.025
$ws.Cells.Item(422, 19).Value = 1.825
$ws.Cells.Item(422, 20).Value = 2
$ws.Cells.Item(422, 21).Value = 1.975
$ws.Cells.Item(422, 22).Value = 1.875
$ws.Cells.Item(422, 23).Value = -1
$ws.Cells.Item(422, 24).Value = -1
$ws.Cells.Item(422, 25).Value = 1.2
$ws.Cells.Item(422, 26).Value = -1
$ws.Cells.Item(422, 27).Value = 0.825
$ws.Cells.Item(422, 28).Value = 0
$ws.Cells.Item(422, 29).Value = -0

# Row 455 <- source row 456
$ws.Cells.Item(455, 2).Value = 5374156
$ws.Cells.Item(455, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(455, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(455, 5).Value = 44857.5625
$ws.Cells.Item(455, 6).Value = 'PAOK Salonika'
$ws.Cells.Item(455, 7).Value = 'Asteras Tripolis'
$ws.Cells.Item(455, 8).Value = 2
$ws.Cells.Item(455, 9).Value = 2
$ws.Cells.Item(455, 10).Value = 'D'
$ws.Cells.Item(455, 11).Value = 1.5
$ws.Cells.Item(455, 12).Value = 4
$ws.Cells.Item(455, 13).Value = 7.5
$ws.Cells.Item(455, 14).Value = 1.363
$ws.Cells.Item(455, 15).Value = 4.333
$ws.Cells.Item(455, 16).Value = 11
$ws.Cells.Item(455, 17).Value = -1.25
$ws.Cells.Item(455, 18).Value = 1.9
$ws.Cells.Item(455, 19).Value = 1.95
$ws.Cells.Item(455, 20).Value = 2.25
$ws.Cells.Item(455, 21).Value = 1.875
$ws.Cells.Item(455, 22).Value = 1.975
$ws.Cells.Item(455, 23).Value = -1
$ws.Cells.Item(455, 24).Value = 3.333
$ws.Cells.Item(455, 25).Value = -1
$ws.Cells.Item(455, 26).Value = -1
$ws.Cells.Item(455, 27).Value = 0.95
$ws.Cells.Item(455, 28).Value = 0.875
$ws.Cells.Item(455, 29).Value = -1

# Row 456 <- source row 455
$ws.Cells.Item(456, 2).Value = 5374154
$ws.Cells.Item(456, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(456, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(456, 5).Value = 44857.5625
$ws.Cells.Item(456, 6).Value = 'Panathinaikos'
$ws.Cells.Item(456, 7).Value = 'Aris Salonika'
$ws.Cells.Item(456, 8).Value = 1
$ws.Cells.Item(456, 9).Value = 0
$ws.Cells.Item(456, 10).Value = 'H'
$ws.Cells.Item(456, 11).Value = 1.8
$ws.Cells.Item(456, 12).Value = 3.4
$ws.Cells.Item(456, 13).Value = 4.75
$ws.Cells.Item(456, 14).Value = 1.727
$ws.Cells.Item(456, 15).Value = 3.4
$ws.Cells.Item(456, 16).Value = 5.75
$ws.Cells.Item(456, 17).Value = -0.75
$ws.Cells.Item(456, 18).Value = 2
$ws.Cells.Item(456, 19).Value = 1.85
$ws.Cells.Item(456, 20).Value = 2
$ws.Cells.Item(456, 21).Value = 2.05
$ws.Cells.Item(456, 22).Value = 1.8
$ws.Cells.Item(456, 23).Value = 0.7270000000000001
$ws.Cells.Item(456, 24).Value = -1
$ws.Cells.Item(456, 25).Value = -1
$ws.Cells.Item(456, 26).Value = 0.5
$ws.Cells.Item(456, 27).Value = -0.5
$ws.Cells.Item(456, 28).Value = -1
$ws.Cells.Item(456, 29).Value = 0.8

# Row 571 <- source row 576
$ws.Cells.Item(571, 2).Value = 5369548
$ws.Cells.Item(571, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(571, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(571, 5).Value = 44997.58333333334
$ws.Cells.Item(571, 6).Value = 'Atromitos Athinon'
$ws.Cells.Item(571, 7).Value = 'Panathinaikos'
$ws.Cells.Item(571, 8).Value = 0
$ws.Cells.Item(571, 9).Value = 2
$ws.Cells.Item(571, 10).Value = 'A'
$ws.Cells.Item(571, 11).Value = 5.75
$ws.Cells.Item(571, 12).Value = 3.25
$ws.Cells.Item(571, 13).Value = 1.75
$ws.Cells.Item(571, 14).Value = 7
$ws.Cells.Item(571, 15).Value = 3.6
$ws.Cells.Item(571, 16).Value = 1.571
$ws.Cells.Item(571, 17).Value = 0.75
$ws.Cells.Item(571, 18).Value = 2.05
$ws.Cells.Item(571, 19).Value = 1.75
$ws.Cells.Item(571, 20).Value = 2.25
$ws.Cells.Item(571, 21).Value = 1.95
$ws.Cells.Item(571, 22).Value = 1.9
$ws.Cells.Item(571, 23).Value = -1
$ws.Cells.Item(571, 24).Value = -1
$ws.Cells.Item(571, 25).Value = 0.571
$ws.Cells.Item(571, 26).Value = -1
$ws.Cells.Item(571, 27).Value = 0.75
$ws.Cells.Item(571, 28).Value = -0.5
$ws.Cells.Item(571, 29).Value = 0.45

# Row 572 <- source row 575
$ws.Cells.Item(572, 2).Value = 5369519
$ws.Cells.Item(572, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(572, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(572, 5).Value = 44997.58333333334
$ws.Cells.Item(572, 6).Value = 'Levadiakos'
$ws.Cells.Item(572, 7).Value = 'OFI Crete'
$ws.Cells.Item(572, 8).Value = 2
$ws.Cells.Item(572, 9).Value = 0
$ws.Cells.Item(572, 10).Value = 'H'
$ws.Cells.Item(572, 11).Value = 2.55
$ws.Cells.Item(572, 12).Value = 3.25
$ws.Cells.Item(572, 13).Value = 2.8
$ws.Cells.Item(572, 14).Value = 2.9
$ws.Cells.Item(572, 15).Value = 3.3
$ws.Cells.Item(572, 16).Value = 2.45
$ws.Cells.Item(572, 17).Value = 0
$ws.Cells.Item(572, 18).Value = 2.075
$ws.Cells.Item(572, 19).Value = 1.725
$ws.Cells.Item(572, 20).Value = 2
$ws.Cells.Item(572, 21).Value = 1.95
$ws.Cells.Item(572, 22).Value = 1.9
$ws.Cells.Item(572, 23).Value = 1.9
$ws.Cells.Item(572, 24).Value = -1
$ws.Cells.Item(572, 25).Value = -1
$ws.Cells.Item(572, 26).Value = 1.075
$ws.Cells.Item(572, 27).Value = -1
$ws.Cells.Item(572, 28).Value = 0
$ws.Cells.Item(572, 29).Value = -0

# Row 573 <- source row 574
$ws.Cells.Item(573, 2).Value = 5374244
$ws.Cells.Item(573, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(573, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(573, 5).Value = 44997.58333333334
$ws.Cells.Item(573, 6).Value = 'Panetolikos'
$ws.Cells.Item(573, 7).Value = 'Lamia'
$ws.Cells.Item(573, 8).Value = 1
$ws.Cells.Item(573, 9).Value = 1
$ws.Cells.Item(573, 10).Value = 'D'
$ws.Cells.Item(573, 11).Value = 2.3
$ws.Cells.Item(573, 12).Value = 3.1
$ws.Cells.Item(573, 13).Value = 3.4
$ws.Cells.Item(573, 14).Value = 2.25
$ws.Cells.Item(573, 15).Value = 3.1
$ws.Cells.Item(573, 16).Value = 3.5
$ws.Cells.Item(573, 17).Value = -0.25
$ws.Cells.Item(573, 18).Value = 1.95
$ws.Cells.Item(573, 19).Value = 1.9
$ws.Cells.Item(573, 20).Value = 2
$ws.Cells.Item(573, 21).Value = 1.95
$ws.Cells.Item(573, 22).Value = 1.9
$ws.Cells.Item(573, 23).Value = -1
$ws.Cells.Item(573, 24).Value = 2.1
$ws.Cells.Item(573, 25).Value = -1
$ws.Cells.Item(573, 26).Value = -0.5
$ws.Cells.Item(573, 27).Value = 0.45
$ws.Cells.Item(573, 28).Value = 0
$ws.Cells.Item(573, 29).Value = -0

# Row 574 <- source row 573
$ws.Cells.Item(574, 2).Value = 5374242
$ws.Cells.Item(574, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(574, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(574, 5).Value = 44997.58333333334
$ws.Cells.Item(574, 6).Value = 'Volos NFC'
$ws.Cells.Item(574, 7).Value = 'PAOK Salonika'
$ws.Cells.Item(574, 8).Value = 0
$ws.Cells.Item(574, 9).Value = 1
$ws.Cells.Item(574, 10).Value = 'A'
$ws.Cells.Item(574, 11).Value = 12
$ws.Cells.Item(574, 12).Value = 5
$ws.Cells.Item(574, 13).Value = 1.3
$ws.Cells.Item(574, 14).Value = 11
$ws.Cells.Item(574, 15).Value = 4.75
$ws.Cells.Item(574, 16).Value = 1.285
$ws.Cells.Item(574, 17).Value = 1.5
$ws.Cells.Item(574, 18).Value = 1.875
$ws.Cells.Item(574, 19).Value = 1.975
$ws.Cells.Item(574, 20).Value = 2.5
$ws.Cells.Item(574, 21).Value = 1.875
$ws.Cells.Item(574, 22).Value = 1.975
$ws.Cells.Item(574, 23).Value = -1
$ws.Cells.Item(574, 24).Value = -1
$ws.Cells.Item(574, 25).Value = 0.2849999999999999
$ws.Cells.Item(574, 26).Value = 0.875
$ws.Cells.Item(574, 27).Value = -1
$ws.Cells.Item(574, 28).Value = -1
$ws.Cells.Item(574, 29).Value = 0.9750000000000001

# Row 575 <- source row 571
$ws.Cells.Item(575, 2).Value = 5374241
$ws.Cells.Item(575, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(575, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(575, 5).Value = 44997.58333333334
$ws.Cells.Item(575, 6).Value = 'Aris Salonika'
$ws.Cells.Item(575, 7).Value = 'Giannina'
$ws.Cells.Item(575, 8).Value = 3
$ws.Cells.Item(575, 9).Value = 1
$ws.Cells.Item(575, 10).Value = 'H'
$ws.Cells.Item(575, 11).Value = 1.4
$ws.Cells.Item(575, 12).Value = 4.333
$ws.Cells.Item(575, 13).Value = 9
$ws.Cells.Item(575, 14).Value = 1.444
$ws.Cells.Item(575, 15).Value = 4
$ws.Cells.Item(575, 16).Value = 9
$ws.Cells.Item(575, 17).Value = -1.25
$ws.Cells.Item(575, 18).Value = 1.975
$ws.Cells.Item(575, 19).Value = 1.875
$ws.Cells.Item(575, 20).Value = 2.25
$ws.Cells.Item(575, 21).Value = 1.85
$ws.Cells.Item(575, 22).Value = 2
$ws.Cells.Item(575, 23).Value = 0.444
$ws.Cells.Item(575, 24).Value = -1
$ws.Cells.Item(575, 25).Value = -1
$ws.Cells.Item(575, 26).Value = 0.9750000000000001
$ws.Cells.Item(575, 27).Value = -1
$ws.Cells.Item(575, 28).Value = 0.8500000000000001
$ws.Cells.Item(575, 29).Value = -1

# Row 576 <- source row 572
$ws.Cells.Item(576, 2).Value = 5374240
$ws.Cells.Item(576, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(576, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(576, 5).Value = 44997.58333333334
$ws.Cells.Item(576, 6).Value = 'AEK Athens'
$ws.Cells.Item(576, 7).Value = 'Olympiakos'
$ws.Cells.Item(576, 8).Value = 1
$ws.Cells.Item(576, 9).Value = 3
$ws.Cells.Item(576, 10).Value = 'A'
$ws.Cells.Item(576, 11).Value = 2.2
$ws.Cells.Item(576, 12).Value = 3.25
$ws.Cells.Item(576, 13).Value = 3.4
$ws.Cells.Item(576, 14).Value = 1.85
$ws.Cells.Item(576, 15).Value = 3.4
$ws.Cells.Item(576, 16).Value = 4.333
$ws.Cells.Item(576, 17).Value = -0.5
$ws.Cells.Item(576, 18).Value = 1.875
$ws.Cells.Item(576, 19).Value = 1.975
$ws.Cells.Item(576, 20).Value = 2.5
$ws.Cells.Item(576, 21).Value = 2
$ws.Cells.Item(576, 22).Value = 1.85
$ws.Cells.Item(576, 23).Value = -1
$ws.Cells.Item(576, 24).Value = -1
$ws.Cells.Item(576, 25).Value = 3.333
$ws.Cells.Item(576, 26).Value = -1
$ws.Cells.Item(576, 27).Value = 0.9750000000000001
$ws.Cells.Item(576, 28).Value = 1
$ws.Cells.Item(576, 29).Value = -1

# Row 604 <- source row 605
$ws.Cells.Item(604, 2).Value = 6399628
$ws.Cells.Item(604, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(604, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(604, 5).Value = 45038.55208333334
$ws.Cells.Item(604, 6).Value = 'Lamia'
$ws.Cells.Item(604, 7).Value = 'Atromitos Athinon'
$ws.Cells.Item(604, 8).Value = 1
$ws.Cells.Item(604, 9).Value = 0
$ws.Cells.Item(604, 10).Value = 'H'
$ws.Cells.Item(604, 11).Value = 2.2
$ws.Cells.Item(604, 12).Value = 3.2
$ws.Cells.Item(604, 13).Value = 3.5
$ws.Cells.Item(604, 14).Value = 1.8
$ws.Cells.Item(604, 15).Value = 3.5
$ws.Cells.Item(604, 16).Value = 4.75
$ws.Cells.Item(604, 17).Value = -0.5
$ws.Cells.Item(604, 18).Value = 1.825
$ws.Cells.Item(604, 19).Value = 2.025
$ws.Cells.Item(604, 20).Value = 2.25
$ws.Cells.Item(604, 21).Value = 2.05
$ws.Cells.Item(604, 22).Value = 1.8
$ws.Cells.Item(604, 23).Value = 0.8
$ws.Cells.Item(604, 24).Value = -1
$ws.Cells.Item(604, 25).Value = -1
$ws.Cells.Item(604, 26).Value = 0.825
$ws.Cells.Item(604, 27).Value = -1
$ws.Cells.Item(604, 28).Value = -1
$ws.Cells.Item(604, 29).Value = 0.8

# Row 605 <- source row 604
$ws.Cells.Item(605, 2).Value = 6399629
$ws.Cells.Item(605, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(605, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(605, 5).Value = 45038.55208333334
$ws.Cells.Item(605, 6).Value = 'Levadiakos'
$ws.Cells.Item(605, 7).Value = 'Ionikos Nikea'
$ws.Cells.Item(605, 8).Value = 2
$ws.Cells.Item(605, 9).Value = 2
$ws.Cells.Item(605, 10).Value = 'D'
$ws.Cells.Item(605, 11).Value = 2.4
$ws.Cells.Item(605, 12).Value = 3.25
$ws.Cells.Item(605, 13).Value = 3
$ws.Cells.Item(605, 14).Value = 2.25
$ws.Cells.Item(605, 15).Value = 3.1
$ws.Cells.Item(605, 16).Value = 3.5
$ws.Cells.Item(605, 17).Value = -0.25
$ws.Cells.Item(605, 18).Value = 1.9
$ws.Cells.Item(605, 19).Value = 1.95
$ws.Cells.Item(605, 20).Value = 1.75
$ws.Cells.Item(605, 21).Value = 1.825
$ws.Cells.Item(605, 22).Value = 2.025
$ws.Cells.Item(605, 23).Value = -1
$ws.Cells.Item(605, 24).Value = 2.1
$ws.Cells.Item(605, 25).Value = -1
$ws.Cells.Item(605, 26).Value = -0.5
$ws.Cells.Item(605, 27).Value = 0.475
$ws.Cells.Item(605, 28).Value = 0.825
$ws.Cells.Item(605, 29).Value = -1

# Row 617 <- source row 618
$ws.Cells.Item(617, 2).Value = 6397967
$ws.Cells.Item(617, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(617, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(617, 5).Value = 45046.58333333334
$ws.Cells.Item(617, 6).Value = 'Panathinaikos'
$ws.Cells.Item(617, 7).Value = 'AEK Athens'
$ws.Cells.Item(617, 8).Value = 0
$ws.Cells.Item(617, 9).Value = 0
$ws.Cells.Item(617, 10).Value = 'D'
$ws.Cells.Item(617, 11).Value = 2.625
$ws.Cells.Item(617, 12).Value = 3
$ws.Cells.Item(617, 13).Value = 2.9
$ws.Cells.Item(617, 14).Value = 2.7
$ws.Cells.Item(617, 15).Value = 3.1
$ws.Cells.Item(617, 16).Value = 2.8
$ws.Cells.Item(617, 17).Value = 0
$ws.Cells.Item(617, 18).Value = 1.9
$ws.Cells.Item(617, 19).Value = 1.95
$ws.Cells.Item(617, 20).Value = 2
$ws.Cells.Item(617, 21).Value = 1.875
$ws.Cells.Item(617, 22).Value = 1.975
$ws.Cells.Item(617, 23).Value = -1
$ws.Cells.Item(617, 24).Value = 2.1
$ws.Cells.Item(617, 25).Value = -1
$ws.Cells.Item(617, 26).Value = 0
$ws.Cells.Item(617, 27).Value = -0
$ws.Cells.Item(617, 28).Value = -1
$ws.Cells.Item(617, 29).Value = 0.9750000000000001

# Row 618 <- source row 617
$ws.Cells.Item(618, 2).Value = 6399613
$ws.Cells.Item(618, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(618, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(618, 5).Value = 45046.58333333334
$ws.Cells.Item(618, 6).Value = 'PAOK Salonika'
$ws.Cells.Item(618, 7).Value = 'Aris Salonika'
$ws.Cells.Item(618, 8).Value = 3
$ws.Cells.Item(618, 9).Value = 2
$ws.Cells.Item(618, 10).Value = 'H'
$ws.Cells.Item(618, 11).Value = 1.75
$ws.Cells.Item(618, 12).Value = 3.4
$ws.Cells.Item(618, 13).Value = 5.25
$ws.Cells.Item(618, 14).Value = 1.7
$ws.Cells.Item(618, 15).Value = 3.6
$ws.Cells.Item(618, 16).Value = 5.25
$ws.Cells.Item(618, 17).Value = -0.75
$ws.Cells.Item(618, 18).Value = 1.875
$ws.Cells.Item(618, 19).Value = 1.975
$ws.Cells.Item(618, 20).Value = 2.25
$ws.Cells.Item(618, 21).Value = 1.8
$ws.Cells.Item(618, 22).Value = 2.05
$ws.Cells.Item(618, 23).Value = 0.7
$ws.Cells.Item(618, 24).Value = -1
$ws.Cells.Item(618, 25).Value = -1
$ws.Cells.Item(618, 26).Value = 0.4375
$ws.Cells.Item(618, 27).Value = -0.5
$ws.Cells.Item(618, 28).Value = 0.8
$ws.Cells.Item(618, 29).Value = -1

# Row 622 <- source row 625
$ws.Cells.Item(622, 2).Value = 6399635
$ws.Cells.Item(622, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(622, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(622, 5).Value = 45052.58333333334
$ws.Cells.Item(622, 6).Value = 'Atromitos Athinon'
$ws.Cells.Item(622, 7).Value = 'Panetolikos'
$ws.Cells.Item(622, 8).Value = 2
$ws.Cells.Item(622, 9).Value = 0
$ws.Cells.Item(622, 10).Value = 'H'
$ws.Cells.Item(622, 11).Value = 2.25
$ws.Cells.Item(622, 12).Value = 3.1
$ws.Cells.Item(622, 13).Value = 3.5
$ws.Cells.Item(622, 14).Value = 1.95
$ws.Cells.Item(622, 15).Value = 3.3
$ws.Cells.Item(622, 16).Value = 4.333
$ws.Cells.Item(622, 17).Value = -0.5
$ws.Cells.Item(622, 18).Value = 2
$ws.Cells.Item(622, 19).Value = 1.85
$ws.Cells.Item(622, 20).Value = 2.25
$ws.Cells.Item(622, 21).Value = 1.925
$ws.Cells.Item(622, 22).Value = 1.925
$ws.Cells.Item(622, 23).Value = 0.95
$ws.Cells.Item(622, 24).Value = -1
$ws.Cells.Item(622, 25).Value = -1
$ws.Cells.Item(622, 26).Value = 1
$ws.Cells.Item(622, 27).Value = -1
$ws.Cells.Item(622, 28).Value = -0.5
$ws.Cells.Item(622, 29).Value = 0.4625

# Row 623 <- source row 624
$ws.Cells.Item(623, 2).Value = 6399638
$ws.Cells.Item(623, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(623, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(623, 5).Value = 45052.58333333334
$ws.Cells.Item(623, 6).Value = 'OFI Crete'
$ws.Cells.Item(623, 7).Value = 'Ionikos Nikea'
$ws.Cells.Item(623, 8).Value = 2
$ws.Cells.Item(623, 9).Value = 2
$ws.Cells.Item(623, 10).Value = 'D'
$ws.Cells.Item(623, 11).Value = 2.55
$ws.Cells.Item(623, 12).Value = 3.25
$ws.Cells.Item(623, 13).Value = 2.8
$ws.Cells.Item(623, 14).Value = 2.45
$ws.Cells.Item(623, 15).Value = 3.1
$ws.Cells.Item(623, 16).Value = 3
$ws.Cells.Item(623, 17).Value = -0.25
$ws.Cells.Item(623, 18).Value = 2.075
$ws.Cells.Item(623, 19).Value = 1.725
$ws.Cells.Item(623, 20).Value = 2.25
$ws.Cells.Item(623, 21).Value = 1.875
$ws.Cells.Item(623, 22).Value = 1.975
$ws.Cells.Item(623, 23).Value = -1
$ws.Cells.Item(623, 24).Value = 2.1
$ws.Cells.Item(623, 25).Value = -1
$ws.Cells.Item(623, 26).Value = -0.5
$ws.Cells.Item(623, 27).Value = 0.3625
$ws.Cells.Item(623, 28).Value = 0.875
$ws.Cells.Item(623, 29).Value = -1

# Row 624 <- source row 623
$ws.Cells.Item(624, 2).Value = 6399637
$ws.Cells.Item(624, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(624, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(624, 5).Value = 45052.58333333334
$ws.Cells.Item(624, 6).Value = 'Lamia'
$ws.Cells.Item(624, 7).Value = 'Levadiakos'
$ws.Cells.Item(624, 8).Value = 1
$ws.Cells.Item(624, 9).Value = 1
$ws.Cells.Item(624, 10).Value = 'D'
$ws.Cells.Item(624, 11).Value = 2.1
$ws.Cells.Item(624, 12).Value = 3.2
$ws.Cells.Item(624, 13).Value = 3.8
$ws.Cells.Item(624, 14).Value = 2.1
$ws.Cells.Item(624, 15).Value = 3.1
$ws.Cells.Item(624, 16).Value = 3.8
$ws.Cells.Item(624, 17).Value = -0.25
$ws.Cells.Item(624, 18).Value = 1.8
$ws.Cells.Item(624, 19).Value = 2.05
$ws.Cells.Item(624, 20).Value = 2
$ws.Cells.Item(624, 21).Value = 2.1
$ws.Cells.Item(624, 22).Value = 1.775
$ws.Cells.Item(624, 23).Value = -1
$ws.Cells.Item(624, 24).Value = 2.1
$ws.Cells.Item(624, 25).Value = -1
$ws.Cells.Item(624, 26).Value = -0.5
$ws.Cells.Item(624, 27).Value = 0.5249999999999999
$ws.Cells.Item(624, 28).Value = 0
$ws.Cells.Item(624, 29).Value = -0

# Row 625 <- source row 622
$ws.Cells.Item(625, 2).Value = 6399636
$ws.Cells.Item(625, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(625, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(625, 5).Value = 45052.58333333334
$ws.Cells.Item(625, 6).Value = 'Giannina'
$ws.Cells.Item(625, 7).Value = 'Asteras Tripolis'
$ws.Cells.Item(625, 8).Value = 1
$ws.Cells.Item(625, 9).Value = 0
$ws.Cells.Item(625, 10).Value = 'H'
$ws.Cells.Item(625, 11).Value = 2.5
$ws.Cells.Item(625, 12).Value = 2.45
$ws.Cells.Item(625, 13).Value = 4.1
$ws.Cells.Item(625, 14).Value = 3
$ws.Cells.Item(625, 15).Value = 2.05
$ws.Cells.Item(625, 16).Value = 4.2
$ws.Cells.Item(625, 17).Value = 0
$ws.Cells.Item(625, 18).Value = 1.875
$ws.Cells.Item(625, 19).Value = 1.975
$ws.Cells.Item(625, 20).Value = 1.75
$ws.Cells.Item(625, 21).Value = 2.025
$ws.Cells.Item(625, 22).Value = 1.825
$ws.Cells.Item(625, 23).Value = 2
$ws.Cells.Item(625, 24).Value = -1
$ws.Cells.Item(625, 25).Value = -1
$ws.Cells.Item(625, 26).Value = 0.875
$ws.Cells.Item(625, 27).Value = -1
$ws.Cells.Item(625, 28).Value = -1
$ws.Cells.Item(625, 29).Value = 0.825

# Row 627 <- source row 628
$ws.Cells.Item(627, 2).Value = 6399615
$ws.Cells.Item(627, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(627, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(627, 5).Value = 45054.625
$ws.Cells.Item(627, 6).Value = 'Aris Salonika'
$ws.Cells.Item(627, 7).Value = 'AEK Athens'
$ws.Cells.Item(627, 8).Value = 1
$ws.Cells.Item(627, 9).Value = 2
$ws.Cells.Item(627, 10).Value = 'A'
$ws.Cells.Item(627, 11).Value = 4.5
$ws.Cells.Item(627, 12).Value = 3.25
$ws.Cells.Item(627, 13).Value = 1.8
$ws.Cells.Item(627, 14).Value = 12
$ws.Cells.Item(627, 15).Value = 5.25
$ws.Cells.Item(627, 16).Value = 1.285
$ws.Cells.Item(627, 17).Value = 1.5
$ws.Cells.Item(627, 18).Value = 2
$ws.Cells.Item(627, 19).Value = 1.85
$ws.Cells.Item(627, 20).Value = 2.75
$ws.Cells.Item(627, 21).Value = 1.875
$ws.Cells.Item(627, 22).Value = 1.975
$ws.Cells.Item(627, 23).Value = -1
$ws.Cells.Item(627, 24).Value = -1
$ws.Cells.Item(627, 25).Value = 0.2849999999999999
$ws.Cells.Item(627, 26).Value = 1
$ws.Cells.Item(627, 27).Value = -1
$ws.Cells.Item(627, 28).Value = 0.4375
$ws.Cells.Item(627, 29).Value = -0.5

# Row 628 <- source row 627
$ws.Cells.Item(628, 2).Value = 6399616
$ws.Cells.Item(628, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(628, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(628, 5).Value = 45054.625
$ws.Cells.Item(628, 6).Value = 'Volos NFC'
$ws.Cells.Item(628, 7).Value = 'PAOK Salonika'
$ws.Cells.Item(628, 8).Value = 0
$ws.Cells.Item(628, 9).Value = 2
$ws.Cells.Item(628, 10).Value = 'A'
$ws.Cells.Item(628, 11).Value = 13
$ws.Cells.Item(628, 12).Value = 6
$ws.Cells.Item(628, 13).Value = 1.181
$ws.Cells.Item(628, 14).Value = 17
$ws.Cells.Item(628, 15).Value = 6
$ws.Cells.Item(628, 16).Value = 1.2
$ws.Cells.Item(628, 17).Value = 2
$ws.Cells.Item(628, 18).Value = 1.85
$ws.Cells.Item(628, 19).Value = 2
$ws.Cells.Item(628, 20).Value = 2.75
$ws.Cells.Item(628, 21).Value = 1.8
$ws.Cells.Item(628, 22).Value = 2.05
$ws.Cells.Item(628, 23).Value = -1
$ws.Cells.Item(628, 24).Value = -1
$ws.Cells.Item(628, 25).Value = 0.2
$ws.Cells.Item(628, 26).Value = 0
$ws.Cells.Item(628, 27).Value = -0
$ws.Cells.Item(628, 28).Value = -1
$ws.Cells.Item(628, 29).Value = 1.05

# Row 630 <- source row 632
$ws.Cells.Item(630, 2).Value = 6399641
$ws.Cells.Item(630, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(630, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(630, 5).Value = 45059.58333333334
$ws.Cells.Item(630, 6).Value = 'Levadiakos'
$ws.Cells.Item(630, 7).Value = 'Giannina'
$ws.Cells.Item(630, 8).Value = 3
$ws.Cells.Item(630, 9).Value = 3
$ws.Cells.Item(630, 10).Value = 'D'
$ws.Cells.Item(630, 11).Value = 2.5
$ws.Cells.Item(630, 12).Value = 3.25
$ws.Cells.Item(630, 13).Value = 2.9
$ws.Cells.Item(630, 14).Value = 2.2
$ws.Cells.Item(630, 15).Value = 3.3
$ws.Cells.Item(630, 16).Value = 3.4
$ws.Cells.Item(630, 17).Value = -0.25
$ws.Cells.Item(630, 18).Value = 1.875
$ws.Cells.Item(630, 19).Value = 1.975
$ws.Cells.Item(630, 20).Value = 2
$ws.Cells.Item(630, 21).Value = 1.925
$ws.Cells.Item(630, 22).Value = 1.925
$ws.Cells.Item(630, 23).Value = -1
$ws.Cells.Item(630, 24).Value = 2.3
$ws.Cells.Item(630, 25).Value = -1
$ws.Cells.Item(630, 26).Value = -0.5
$ws.Cells.Item(630, 27).Value = 0.4875
$ws.Cells.Item(630, 28).Value = 0.925
$ws.Cells.Item(630, 29).Value = -1

# Row 631 <- source row 630
$ws.Cells.Item(631, 2).Value = 6399640
$ws.Cells.Item(631, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(631, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(631, 5).Value = 45059.58333333334
$ws.Cells.Item(631, 6).Value = 'Ionikos Nikea'
$ws.Cells.Item(631, 7).Value = 'Lamia'
$ws.Cells.Item(631, 8).Value = 2
$ws.Cells.Item(631, 9).Value = 2
$ws.Cells.Item(631, 10).Value = 'D'
$ws.Cells.Item(631, 11).Value = 2.15
$ws.Cells.Item(631, 12).Value = 3.25
$ws.Cells.Item(631, 13).Value = 3.6
$ws.Cells.Item(631, 14).Value = 2.1
$ws.Cells.Item(631, 15).Value = 3.1
$ws.Cells.Item(631, 16).Value = 3.8
$ws.Cells.Item(631, 17).Value = -0.25
$ws.Cells.Item(631, 18).Value = 1.825
$ws.Cells.Item(631, 19).Value = 2.025
$ws.Cells.Item(631, 20).Value = 2
$ws.Cells.Item(631, 21).Value = 1.8
$ws.Cells.Item(631, 22).Value = 2.05
$ws.Cells.Item(631, 23).Value = -1
$ws.Cells.Item(631, 24).Value = 2.1
$ws.Cells.Item(631, 25).Value = -1
$ws.Cells.Item(631, 26).Value = -0.5
$ws.Cells.Item(631, 27).Value = 0.5125
$ws.Cells.Item(631, 28).Value = 0.8
$ws.Cells.Item(631, 29).Value = -1

# Row 632 <- source row 631
$ws.Cells.Item(632, 2).Value = 6399639
$ws.Cells.Item(632, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(632, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(632, 5).Value = 45059.58333333334
$ws.Cells.Item(632, 6).Value = 'Asteras Tripolis'
$ws.Cells.Item(632, 7).Value = 'Atromitos Athinon'
$ws.Cells.Item(632, 8).Value = 1
$ws.Cells.Item(632, 9).Value = 1
$ws.Cells.Item(632, 10).Value = 'D'
$ws.Cells.Item(632, 11).Value = 2.35
$ws.Cells.Item(632, 12).Value = 3.25
$ws.Cells.Item(632, 13).Value = 3.1
$ws.Cells.Item(632, 14).Value = 2.7
$ws.Cells.Item(632, 15).Value = 3.3
$ws.Cells.Item(632, 16).Value = 2.6
$ws.Cells.Item(632, 17).Value = 0
$ws.Cells.Item(632, 18).Value = 1.925
$ws.Cells.Item(632, 19).Value = 1.925
$ws.Cells.Item(632, 20).Value = 2.25
$ws.Cells.Item(632, 21).Value = 2.05
$ws.Cells.Item(632, 22).Value = 1.8
$ws.Cells.Item(632, 23).Value = -1
$ws.Cells.Item(632, 24).Value = 2.3
$ws.Cells.Item(632, 25).Value = -1
$ws.Cells.Item(632, 26).Value = 0
$ws.Cells.Item(632, 27).Value = -0
$ws.Cells.Item(632, 28).Value = -0.5
$ws.Cells.Item(632, 29).Value = 0.4

# Row 685 <- source row 686
$ws.Cells.Item(685, 2).Value = 6937192
$ws.Cells.Item(685, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(685, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(685, 5).Value = 45221.41666666666
$ws.Cells.Item(685, 6).Value = 'Asteras Tripolis'
$ws.Cells.Item(685, 7).Value = 'AEK Athens'
$ws.Cells.Item(685, 8).Value = 0
$ws.Cells.Item(685, 9).Value = 3
$ws.Cells.Item(685, 10).Value = 'A'
$ws.Cells.Item(685, 11).Value = 6.5
$ws.Cells.Item(685, 12).Value = 3.6
$ws.Cells.Item(685, 13).Value = 1.6
$ws.Cells.Item(685, 14).Value = 7.5
$ws.Cells.Item(685, 15).Value = 4
$ws.Cells.Item(685, 16).Value = 1.5
$ws.Cells.Item(685, 17).Value = 1
$ws.Cells.Item(685, 18).Value = 2.025
$ws.Cells.Item(685, 19).Value = 1.825
$ws.Cells.Item(685, 20).Value = 2.5
$ws.Cells.Item(685, 21).Value = 1.925
$ws.Cells.Item(685, 22).Value = 1.925
$ws.Cells.Item(685, 23).Value = -1
$ws.Cells.Item(685, 24).Value = -1
$ws.Cells.Item(685, 25).Value = 0.5
$ws.Cells.Item(685, 26).Value = -1
$ws.Cells.Item(685, 27).Value = 0.825
$ws.Cells.Item(685, 28).Value = 0.925
$ws.Cells.Item(685, 29).Value = -1

# Row 686 <- source row 685
$ws.Cells.Item(686, 2).Value = 6935736
$ws.Cells.Item(686, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(686, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(686, 5).Value = 45221.41666666666
$ws.Cells.Item(686, 6).Value = 'Kifisias FC'
$ws.Cells.Item(686, 7).Value = 'OFI Crete'
$ws.Cells.Item(686, 8).Value = 0
$ws.Cells.Item(686, 9).Value = 0
$ws.Cells.Item(686, 10).Value = 'D'
$ws.Cells.Item(686, 11).Value = 3.3
$ws.Cells.Item(686, 12).Value = 3.25
$ws.Cells.Item(686, 13).Value = 2.25
$ws.Cells.Item(686, 14).Value = 3.1
$ws.Cells.Item(686, 15).Value = 3.4
$ws.Cells.Item(686, 16).Value = 2.4
$ws.Cells.Item(686, 17).Value = 0.25
$ws.Cells.Item(686, 18).Value = 1.8
$ws.Cells.Item(686, 19).Value = 2.05
$ws.Cells.Item(686, 20).Value = 2.5
$ws.Cells.Item(686, 21).Value = 1.825
$ws.Cells.Item(686, 22).Value = 2.025
$ws.Cells.Item(686, 23).Value = -1
$ws.Cells.Item(686, 24).Value = 2.4
$ws.Cells.Item(686, 25).Value = -1
$ws.Cells.Item(686, 26).Value = 0.4
$ws.Cells.Item(686, 27).Value = -0.5
$ws.Cells.Item(686, 28).Value = -1
$ws.Cells.Item(686, 29).Value = 1.025

# Row 758 <- source row 759
$ws.Cells.Item(758, 2).Value = 6937238
$ws.Cells.Item(758, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(758, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(758, 5).Value = 45305.64583333334
$ws.Cells.Item(758, 6).Value = 'PAOK Salonika'
$ws.Cells.Item(758, 7).Value = 'Giannina'
$ws.Cells.Item(758, 8).Value = 4
$ws.Cells.Item(758, 9).Value = 0
$ws.Cells.Item(758, 10).Value = 'H'
$ws.Cells.Item(758, 11).Value = 1.111
$ws.Cells.Item(758, 12).Value = 9
$ws.Cells.Item(758, 13).Value = 23
$ws.Cells.Item(758, 14).Value = 1.25
$ws.Cells.Item(758, 15).Value = 6
$ws.Cells.Item(758, 16).Value = 9
$ws.Cells.Item(758, 17).Value = -1.75
$ws.Cells.Item(758, 18).Value = 2.025
$ws.Cells.Item(758, 19).Value = 1.825
$ws.Cells.Item(758, 20).Value = 2.75
$ws.Cells.Item(758, 21).Value = 1.8
$ws.Cells.Item(758, 22).Value = 2.05
$ws.Cells.Item(758, 23).Value = 0.25
$ws.Cells.Item(758, 24).Value = -1
$ws.Cells.Item(758, 25).Value = -1
$ws.Cells.Item(758, 26).Value = 1.025
$ws.Cells.Item(758, 27).Value = -1
$ws.Cells.Item(758, 28).Value = 0.8
$ws.Cells.Item(758, 29).Value = -1

# Row 759 <- source row 758
$ws.Cells.Item(759, 2).Value = 6936857
$ws.Cells.Item(759, 3).Value = 'Greece Super League 1'
$ws.Cells.Item(759, 4).Value = 'Greece Super League 1'
$ws.Cells.Item(759, 5).Value = 45305.64583333334
$ws.Cells.Item(759, 6).Value = 'AEK Athens'
$ws.Cells.Item(759, 7).Value = 'Panathinaikos'
$ws.Cells.Item(759, 8).Value = 2
$ws.Cells.Item(759, 9).Value = 2
$ws.Cells.Item(759, 10).Value = 'D'
$ws.Cells.Item(759, 11).Value = 1.909
$ws.Cells.Item(759, 12).Value = 3.5
$ws.Cells.Item(759, 13).Value = 4.2
$ws.Cells.Item(759, 14).Value = 2.15
$ws.Cells.Item(759, 15).Value = 3.2
$ws.Cells.Item(759, 16).Value = 3.5
$ws.Cells.Item(759, 17).Value = -0.25
$ws.Cells.Item(759, 18).Value = 1.85
$ws.Cells.Item(759, 19).Value = 2
$ws.Cells.Item(759, 20).Value = 2
$ws.Cells.Item(759, 21).Value = 1.8
$ws.Cells.Item(759, 22).Value = 2.05
$ws.Cells.Item(759, 23).Value = -1
$ws.Cells.Item(759, 24).Value = 2.2
$ws.Cells.Item(759, 25).Value = -1
$ws.Cells.Item(759, 26).Value = -0.5
$ws.Cells.Item(759, 27).Value = 0.5
$ws.Cells.Item(759, 28).Value = 0.8
$ws.Cells.Item(759, 29).Value = -1
# Rows 775-779: in-place odds updates (no row swap)
# Row 775
$ws.Cells.Item(775, 21).Value = 1.95
$ws.Cells.Item(775, 22).Value = 1.9

# Row 776
$ws.Cells.Item(776, 14).Value = 1.75
$ws.Cells.Item(776, 16).Value = 4.75
$ws.Cells.Item(776, 18).Value = 2
$ws.Cells.Item(776, 19).Value = 1.85
$ws.Cells.Item(776, 21).Value = 1.825
$ws.Cells.Item(776, 22).Value = 2.025

# Row 777
$ws.Cells.Item(777, 18).Value = 1.85
$ws.Cells.Item(777, 19).Value = 2
$ws.Cells.Item(777, 21).Value = 1.85
$ws.Cells.Item(777, 22).Value = 2

# Row 779
$ws.Cells.Item(779, 14).Value = 7
$ws.Cells.Item(779, 15).Value = 4.2
$ws.Cells.Item(779, 16).Value = 1.5
$ws.Cells.Item(779, 17).Value = 1
$ws.Cells.Item(779, 18).Value = 2.05
$ws.Cells.Item(779, 19).Value = 1.8
